# issue #5: add legislator_id, name, date into dataframe
#
# The "股票" (stocks) worksheet (4th sheet) gains three new trailing
# columns: date, legislator_name, legislator_id - populated for every
# data row with the report date, the legislator's name and id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$lastRow = 23

# --- Header row (row 1): new headers in H1:J1, matching the bold/bordered
# style already used by the other headers in row 1 (copy from G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows: same values repeated for every stock holding row.
for ($r = 2; $r -le $lastRow; $r++) {
    # Force text so the date-like string isn't reinterpreted as a date serial.
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value = "2012-04-19"
    $ws.Cells.Item($r, 9).Value = "張慶忠"
    $ws.Cells.Item($r, 10).Value = 1347
}
